# Update burndown chart data/title for "Iteracion 3"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 2 ("Real") values ---
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = 21
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 0

# --- Row 3 ("Estimado") starting value (rest are formula-driven) ---
$ws.Range("B3").Value = 24

# --- Chart title: " Iteracion 1" -> " Iteracion 3" ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Burndown chart Iteracion 3"

# --- Selection moves from G3 to K3 ---
$ws.Range("K3").Select()
